$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row 2, shifting all subsequent rows up by one.
$ws.Rows("2:2").Delete()

# Select row 2 (the entire row) to match the resulting selection state.
$ws.Range("A2:XFD2").Select()
